$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.867.04'
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('D3').Value = '1.978.12'
$ws.Range('E3').Value = '  +0.53%  '
$ws.Range('E4').Value = '  +0.04%  '
$r = $ws.Range('D5')
$r.Value = "'245.29"
$r.Style = 'Normal'
$ws.Range('E5').Value = '  +0.20%  '
$ws.Range('E6').Value = '  +1.34%  '
$r = $ws.Range('D7')
$r.Value = "'60.91"
$r.Style = 'Normal'
$ws.Range('E7').Value = '  +2.72%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +1.97%  '
$r = $ws.Range('D10')
$r.Value = "'0.0803"
$r.Style = 'Normal'
$ws.Range('E10').Value = '  -2.08%  '
$ws.Range('E11').Value = '  +0.69%  '
$r = $ws.Range('D12')
$r.Value = "'14.74"
$r.Style = 'Normal'
$ws.Range('E12').Value = '  +6.73%  '
$r = $ws.Range('D13')
$r.Value = "'0.846"
$r.Style = 'Normal'
$ws.Range('E13').Value = '  +1.74%  '
$r = $ws.Range('D14')
$r.Value = "'21.99"
$r.Style = 'Normal'
$ws.Range('E14').Value = '  -1.91%  '
$ws.Range('D15').Value = '2.268.38'
$ws.Range('E15').Value = '  +0.55%  '
$ws.Range('E16').Value = '  +2.77%  '
$ws.Range('D17').Value = '1.975.10'
$ws.Range('E17').Value = '  +0.30%  '
$ws.Range('D18').Value = '36.776.95'
$ws.Range('E18').Value = '  +0.32%  '
$r = $ws.Range('D19')
$r.Value = "'70.14"
$r.Style = 'Normal'
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('D20').Value = '0.0₃0862'
$ws.Range('E20').Value = '  -0.21%  '
$ws.Range('E21').Value = '  +1.14%  '
$r = $ws.Range('D22')
$r.Value = "'230.12"
$r.Style = 'Normal'
$ws.Range('E22').Value = '  +0.31%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('E24').Value = '  +1.84%  '
$ws.Range('E25').Value = '  -0.15%  '
$r = $ws.Range('D26')
$r.Value = "'0.146"
$r.Style = 'Normal'
$ws.Range('E26').Value = '  +2.36%  '
$r = $ws.Range('D27')
$r.Value = "'9.30"
$r.Style = 'Normal'
$ws.Range('E27').Value = '  -0.56%  '
$r = $ws.Range('D28')
$r.Value = "'163.66"
$r.Style = 'Normal'
$ws.Range('E28').Value = '  +1.77%  '
$r = $ws.Range('D29')
$r.Value = "'19.49"
$r.Style = 'Normal'
$ws.Range('E29').Value = '  +0.24%  '
$ws.Range('E30').Value = '  +19.36%  '
$ws.Range('E31').Value = '  +1.43%  '
$ws.Range('E32').Value = '  +2.49%  '
$ws.Range('E33').Value = '  -0.30%  '
$ws.Range('E34').Value = '  +4.84%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$r = $ws.Range('D35')
$r.Value = "'2.27"
$r.Style = 'Normal'
$ws.Range('E35').Value = '  -0.20%  '
$ws.Range('B36').Value = 'BinanceUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$r = $ws.Range('D36')
$r.Value = "'1.00"
$r.Style = 'Normal'
$ws.Range('E36').Value = '  -0.06%  '
$r = $ws.Range('D37')
$r.Value = "'3.35"
$r.Style = 'Normal'
$ws.Range('E37').Value = '  -1.69%  '
$ws.Range('E38').Value = '  +0.21%  '
$ws.Range('E39').Value = '  -10.10%  '
$r = $ws.Range('D40')
$r.Value = "'0.0981"
$r.Style = 'Normal'
$ws.Range('E40').Value = '  -1.73%  '
$ws.Range('E41').Value = '  +0.67%  '
$ws.Range('E42').Value = '  +0.55%  '
$ws.Range('E43').Value = '  -0.04%  '
$r = $ws.Range('D44')
$r.Value = "'16.26"
$r.Style = 'Normal'
$ws.Range('E44').Value = '  +0.59%  '
$ws.Range('D45').Value = '1.368.54'
$ws.Range('E45').Value = '  +0.20%  '
$r = $ws.Range('D46')
$r.Value = "'89.84"
$r.Style = 'Normal'
$ws.Range('E46').Value = '  +2.01%  '
$ws.Range('E47').Value = '  -0.26%  '
$r = $ws.Range('D48')
$r.Value = "'7.25"
$r.Style = 'Normal'
$ws.Range('E48').Value = '  +0.93%  '
$r = $ws.Range('D49')
$r.Value = "'2.82"
$r.Style = 'Normal'
$ws.Range('E49').Value = '  -0.82%  '
$r = $ws.Range('D50')
$r.Value = "'46.30"
$r.Style = 'Normal'
$ws.Range('E50').Value = '  +4.93%  '
$r = $ws.Range('D51')
$r.Value = "'1.96"
$r.Style = 'Normal'
$ws.Range('E51').Value = '  +9.69%  '
